$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-15.11113567023471, 0.9677864928380381, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 1.091925508942268, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, 0.9640300701901376, -15.11113567023471, -15.11113567023471, 3.785507156932909, -15.11113567023471, 1.818134298412369, -15.11113567023471, 2.569038394384597, -15.11113567023471),
    @(-15.11113567023471, 1.379793529987583, -15.11113567023471, -15.11113567023471, -15.11113567023471, 3.400543642611971, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(2.610477568607593, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, 1.994210198669613, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(3.79622741045806, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 1.275552027482065, -15.11113567023471, 2.106945048173523),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, 3.04061919649519, -15.11113567023471, 1.814893872126074, -15.11113567023471, -15.11113567023471, -15.11113567023471, 1.734366631804398),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, 1.802092693164636, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 1.678410710334815, 2.495479187650318),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 1.779472992017555),
    @(-15.11113567023471, -15.11113567023471, 4.321889369223833, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 2.011904588916694, -15.11113567023471),
    @(-15.11113567023471, 1.614882869827343, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 1.949528600838573, 0.8885340852807125, 1.845614816721936, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 1.925048676616946, 0.4717413527805909, 1.694872638189281, -15.11113567023471),
    @(-15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, -15.11113567023471, 2.06528099883354, 1.542842963394889, -15.11113567023471, -15.11113567023471),
    @(-15.11113567023471, 1.98167307715685, -15.11113567023471, -15.11113567023471, 2.634571475336438, -15.11113567023471, 1.110545845387577, 3.216509069021842, -15.11113567023471, 1.71909479226271),
    @(-15.11113567023471, 2.6956216111858, -15.11113567023471, 2.104293051894625, -15.11113567023471, 2.565875055706133, 1.303896325588812, -15.11113567023471, -15.11113567023471, -15.11113567023471)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item(2 + $r, 2 + $c).Value = $row[$c]
    }
}
